# Automatische test-sync: 2025-07-23 22:54:50
#
# Adds a new log entry (row 28) to the "Logs" sheet, extends the
# column D/G/H/I/J conditional-formatting ranges to cover the new row,
# and swaps the "IT / Technisch probleem" / "Productinformatie" rows
# (and their counts) on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Logs" sheet: append new row 28
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A28").Value = "Kun je nagaan of we nog EcoPro-700 op voorraad hebben?"
$logs.Range("B28").Value = "mailmind.test@zohomail.eu"
$logs.Range("C28").Value = "Testmail #18: Kun je nagaan of we nog EcoPro-700 op voorraad hebben?"
$logs.Range("D28").Value = "Productinformatie"
$logs.Range("E28").Value = "Beste afzender,`nBedankt voor je bericht. Op dit moment hebben we nog EcoPro-700 op voorraad. Als je wilt bestellen, laat het ons dan weten en we helpen je graag verder.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Range("F28").Value = "2025-07-23 22:54:16"
$logs.Range("G28").Value = "Ja"
$logs.Range("H28").Value = "Nee"
$logs.Range("I28").Value = "Ja"
$logs.Range("J28").Value = "Nee"

# Writing the multi-line E28 text causes an implicit custom row height;
# AutoFit puts the row back to the sheet's default (no ht/customHeight
# attributes emitted), matching the untouched rows above it.
$logs.Rows.Item(28).AutoFit()

# ---------------------------------------------------------------------
# 2. Extend conditional formatting ranges from row 27 to row 28
# ---------------------------------------------------------------------
$logs.Range("D2:D27").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D28"))
$logs.Range("G2:G27").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G28"))
$logs.Range("H2:H27").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H28"))
$logs.Range("I2:I27").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I28"))
$logs.Range("J2:J27").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J28"))

# ---------------------------------------------------------------------
# 3. "Dashboard" sheet: category counts changed because of the new row
#    - "IT / Technisch probleem" and "Productinformatie" swap rows
#      (row 6 <-> row 9), and "Productinformatie" count goes 1 -> 2.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A6").Value = "Productinformatie"
$dash.Range("B6").Value = 2

$dash.Range("A9").Value = "IT / Technisch probleem"
$dash.Range("B9").Value = 1
